# KWL_C14_dates.xlsx: rename the "Details" (column G) entries from the
# old "Midden, H0xx" / "Burial/burial, M0xx" style to "Midden xx" / "Burial xx".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Midden 26"
$ws.Range("G8").Value = "Midden 26"
$ws.Range("G9").Value = "Burial 9"
$ws.Range("G12").Value = "Midden 32"
$ws.Range("G14").Value = "Burial 95"

# Leave the selection where the author's last save left it.
$ws.Range("G15").Select()
